$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has data rows through row 13 (two "Papaya" records dated
# 2021-07-23 / 44400 and 2021-05-20 / 44336). This edit appends a new weekly
# report: the existing rows 12-13 become the new "current" week's records
# (updated values), and their previous content is preserved by pushing it
# down into two brand-new rows 14-15.

# 1) Make room: insert two fresh blank rows right after the current last row.
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(15).Insert()

# 2) Preserve the old row 12 / row 13 data by copying it down into the new
#    rows 14 / 15 before we overwrite rows 12/13 with the new week's values.
$ws.Range("A12:T12").Copy()
$ws.Range("A14:T14").PasteSpecial()
$ws.Range("A13:T13").Copy()
$ws.Range("A15:T15").PasteSpecial()
$excel.CutCopyMode = $false

# 3) Update row 12 in place with the new week's "Primera" record.
$ws.Cells.Item(12, 4).Value = 44904
$ws.Cells.Item(12, 13).Value = 45
$ws.Cells.Item(12, 14).Value = 15000
$ws.Cells.Item(12, 15).Value = 15000
$ws.Cells.Item(12, 16).Value = 15000
$ws.Cells.Item(12, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(12, 20).Value = 10

# 4) Update row 13 in place with the new week's "Segunda" record.
$ws.Cells.Item(13, 4).Value = 44904
$ws.Cells.Item(13, 12).Value = "Segunda"
$ws.Cells.Item(13, 13).Value = 60
$ws.Cells.Item(13, 14).Value = 10000
$ws.Cells.Item(13, 15).Value = 10000
$ws.Cells.Item(13, 16).Value = 10000
$ws.Cells.Item(13, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(13, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(13, 19).Value = 1000
$ws.Cells.Item(13, 20).Value = 10
